$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) stores numeric-looking values as text
# (e.g. "26.482.55", "0.000006959", "106.00"). Force each rewritten
# cell to text format first so Excel keeps the literal digits instead
# of reinterpreting them as numbers/dates and dropping formatting
# (trailing zeros, thousands dots, etc.).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.482.55"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.728.28"
$ws.Range("E3").Value = "  +0.64%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.78"
$ws.Range("E5").Value = "  +2.29%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4797"
$ws.Range("E7").Value = "  +1.54%  "
$ws.Range("E8").Value = "  +1.81%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06225"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.727.77"
$ws.Range("E10").Value = "  +1.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07134"
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.74"
$ws.Range("E12").Value = "  +3.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6187"
$ws.Range("E13").Value = "  +4.31%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.530"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.24"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9995"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.494.25"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9996"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006959"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.951.19"
$ws.Range("E21").Value = "  +1.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.538"
$ws.Range("E22").Value = "  -0.60%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.945"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.301"
$ws.Range("E24").Value = "  -1.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.49"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.36"
$ws.Range("E26").Value = "  +1.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.804"
$ws.Range("E27").Value = "  +2.01%  "
$ws.Range("E28").Value = "  +0.16%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.82"
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.983"
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08026"
$ws.Range("E31").Value = "  +3.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.732"
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04569"
$ws.Range("E33").Value = "  +3.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9988"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.616"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6372"
$ws.Range("E36").Value = "  +2.27%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9895"
$ws.Range("E37").Value = "  +1.82%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9365"
$ws.Range("E38").Value = "  +2.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.065"
$ws.Range("E39").Value = "  +8.45%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.413"
$ws.Range("E40").Value = "  -0.35%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.00"
$ws.Range("E41").Value = "  -5.50%  "
$ws.Range("E42").Value = "  +0.59%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.712"
$ws.Range("E43").Value = "  +10.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01500"
$ws.Range("E44").Value = "  +2.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3913"
$ws.Range("E45").Value = "  +2.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.989"
$ws.Range("E46").Value = "  +11.71%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1191"
$ws.Range("E47").Value = "  +3.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05322"
$ws.Range("E48").Value = "  +0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.09"
$ws.Range("E49").Value = "  +1.21%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.897"
$ws.Range("E50").Value = "  +3.01%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.267"
$ws.Range("E51").Value = "  +3.47%  "
